# Natmi following Dr Hou advice
# Replace the LR-pair result table (rows 2-7 originally, now rows 2-10)
# with an updated dataset that adds "ECs" as a third sending cluster and
# recomputed summary statistics for all sending/target cluster pairs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New body data (20 columns: A..T), one row per sending/target cluster pair.
$data = @(
    @("ECs",  "Hgf", "Cd44", "ECs",  2, 0.6666666666666666, 5.928568666666667, 17.785706, 0.3809768389628236, 0.3809768389628236, 3, 1, 297.8183156666666, 893.4549469999999, 0.8852156413092672, 0.8852156413092673, 1765.636334620842, 15890.72701158758, 0.3372466568264533, 0.3372466568264534),
    @("ECs",  "Hgf", "Cd44", "FAPs", 2, 0.6666666666666666, 5.928568666666667, 17.785706, 0.3809768389628236, 0.3809768389628236, 3, 1, 24.34034433333333, 73.021033, 0.07234764413494278, 0.0723476441349428, 144.3034027504776, 1298.730624754298, 0.02756277676893777, 0.02756277676893778),
    @("ECs",  "Hgf", "Cd44", "sCs",  2, 0.6666666666666666, 5.928568666666667, 17.785706, 0.3809768389628236, 0.3809768389628236, 3, 1, 14.277234, 42.831702, 0.04243671455578994, 0.04243671455578994, 84.64356213906801, 761.792059251612, 0.0161674053674325, 0.0161674053674325),
    @("FAPs", "Hgf", "Cd44", "ECs",  3, 1, 8.583520999999999, 25.750563, 0.5515872180307627, 0.5515872180307626, 3, 1, 297.8183156666666, 893.4549469999999, 0.8852156413092672, 0.8852156413092673, 2556.329766709462, 23006.96790038516, 0.4882736329470962, 0.4882736329470962),
    @("FAPs", "Hgf", "Cd44", "FAPs", 3, 1, 8.583520999999999, 25.750563, 0.5515872180307627, 0.5515872180307626, 3, 1, 24.34034433333333, 73.021033, 0.07234764413494278, 0.0723476441349428, 208.9258567323976, 1880.332710591579, 0.03990603575947271, 0.03990603575947271),
    @("FAPs", "Hgf", "Cd44", "sCs",  3, 1, 8.583520999999999, 25.750563, 0.5515872180307627, 0.5515872180307626, 3, 1, 14.277234, 42.831702, 0.04243671455578994, 0.04243671455578994, 122.548937860914, 1102.940440748226, 0.02340754932419374, 0.02340754932419374),
    @("sCs",  "Hgf", "Cd44", "ECs",  3, 1, 1.049404, 3.148212, 0.06743594300641363, 0.06743594300641362, 3, 1, 297.8183156666666, 893.4549469999999, 0.8852156413092672, 0.8852156413092673, 312.5317317338627, 2812.785585604764, 0.05969535153571764, 0.05969535153571763),
    @("sCs",  "Hgf", "Cd44", "FAPs", 3, 1, 1.049404, 3.148212, 0.06743594300641363, 0.06743594300641362, 3, 1, 24.34034433333333, 73.021033, 0.07234764413494278, 0.0723476441349428, 25.54285470477733, 229.885692342996, 0.004878831606532297, 0.004878831606532297),
    @("sCs",  "Hgf", "Cd44", "sCs",  3, 1, 1.049404, 3.148212, 0.06743594300641363, 0.06743594300641362, 3, 1, 14.277234, 42.831702, 0.04243671455578994, 0.04243671455578994, 14.982586468536, 134.843278216824, 0.002861759864163694, 0.002861759864163694)
)

# Clear the old body rows (the original sheet only went down to row 7) and
# write the fresh block covering rows 2..10.
$ws.Range("A2:T20").Clear()

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $row[$c]
    }
}
